$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B31: updated modified timestamp
$ws.Range("B31").Value = "2022-06-04T11:15:15+00:00"

# Row 34: clear J34:R34 (move content down to row 35)
$ws.Range("J34:R34").Value = ""

# Row 35: set J35:Q35 with the iop/puv/sosa properties (R35 already has MFI, keep it)
$ws.Range("J35").Value = "iop:hasProperty"
$ws.Range("K35").Value = "iop:hasObjectOfInterest"
$ws.Range("L35").Value = "iop:hasMatrix"
$ws.Range("M35").Value = 'iop:hasContextObject(separator=",")'
$ws.Range("N35").Value = 'iop:hasConstraint(separator=",")'
$ws.Range("O35").Value = 'puv:statistic(separator=",")'
$ws.Range("P35").Value = 'puv:usesMethod(separator=",")'
$ws.Range("Q35").Value = 'sosa:madeBySensor(separator=",")'

# Row 39: R39 changes from "MFI,%" to "MFI"
$ws.Range("R39").Value = "MFI"

# Row 44: R44 changes from empty to "MFI,%"
$ws.Range("R44").Value = "MFI,%"
